$wb = $excel.ActiveWorkbook

# ---- LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 07:12:47"
$ws.Cells.Item(3, 1).Value = "Total filas: 82"
$ws.Cells.Item(9, 3).Value = "15_ABASTO"
$ws.Cells.Item(10, 3).Value = "215_EL PELIGRO"
$ws.Cells.Item(49, 1).Value = "06:53:56"
$ws.Cells.Item(49, 3).Value = "17_ROMERO"
$ws.Cells.Item(49, 4).Value = 1
$ws.Cells.Item(50, 1).Value = "06:46:37"
$ws.Cells.Item(50, 3).Value = "14_ABASTO"
$ws.Cells.Item(50, 4).Value = 8
$ws.Cells.Item(55, 1).Value = "07:12:47"
$ws.Cells.Item(55, 4).Value = 1
$ws.Cells.Item(57, 1).Value = "07:12:47"
$ws.Cells.Item(57, 4).Value = 8
$ws.Cells.Item(59, 1).Value = "07:12:47"
$ws.Cells.Item(59, 4).Value = 11
$ws.Cells.Item(62, 1).Value = "07:12:47"
$ws.Cells.Item(62, 4).Value = 17
$ws.Cells.Item(63, 1).Value = "07:12:47"
$ws.Cells.Item(63, 4).Value = 21
$ws.Cells.Item(64, 1).Value = "07:12:47"
$ws.Cells.Item(64, 4).Value = 23
$ws.Cells.Item(65, 1).Value = "07:12:47"
$ws.Cells.Item(65, 4).Value = 24
$ws.Cells.Item(68, 1).Value = "07:12:47"
$ws.Cells.Item(68, 2).Value = "07:41"
$ws.Cells.Item(68, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(68, 4).Value = 29
$ws.Cells.Item(69, 1).Value = "07:12:47"
$ws.Cells.Item(69, 2).Value = "07:43"
$ws.Cells.Item(69, 4).Value = 31
$ws.Cells.Item(70, 1).Value = "06:18:01"
$ws.Cells.Item(70, 2).Value = "07:44"
$ws.Cells.Item(70, 3).Value = "10_OLMOS"
$ws.Cells.Item(70, 4).Value = 86
$ws.Cells.Item(71, 1).Value = "07:12:47"
$ws.Cells.Item(71, 2).Value = "07:49"
$ws.Cells.Item(71, 3).Value = "15_ABASTO"
$ws.Cells.Item(71, 4).Value = 37
$ws.Cells.Item(72, 1).Value = "07:12:47"
$ws.Cells.Item(72, 2).Value = "07:58"
$ws.Cells.Item(72, 4).Value = 46
$ws.Cells.Item(73, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(74, 1).Value = "07:12:47"
$ws.Cells.Item(74, 2).Value = "07:59"
$ws.Cells.Item(74, 4).Value = 47
$ws.Cells.Item(75, 1).Value = "06:18:01"
$ws.Cells.Item(75, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(75, 4).Value = 102
$ws.Cells.Item(76, 1).Value = "06:46:37"
$ws.Cells.Item(76, 2).Value = "08:00"
$ws.Cells.Item(76, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(76, 4).Value = 74
$ws.Cells.Item(77, 2).Value = "08:01"
$ws.Cells.Item(77, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(77, 4).Value = 68
$ws.Cells.Item(78, 1).Value = "07:12:47"
$ws.Cells.Item(78, 2).Value = "08:03"
$ws.Cells.Item(78, 3).Value = "17X38_ROMERO"
$ws.Cells.Item(78, 4).Value = 51
$ws.Cells.Item(79, 1).Value = "06:53:56"
$ws.Cells.Item(79, 2).Value = "08:13"
$ws.Cells.Item(79, 4).Value = 80
$ws.Cells.Item(80, 1).Value = "07:12:47"
$ws.Cells.Item(80, 2).Value = "08:14"
$ws.Cells.Item(80, 3).Value = "10_OLMOS"
$ws.Cells.Item(80, 4).Value = 62
$ws.Cells.Item(81, 1).Value = "07:12:47"
$ws.Cells.Item(81, 2).Value = "08:19"
$ws.Cells.Item(81, 3).Value = "17_ROMERO"
$ws.Cells.Item(81, 4).Value = 67
$ws.Cells.Item(82, 1).Value = "07:12:47"
$ws.Cells.Item(82, 2).Value = "08:29"
$ws.Cells.Item(82, 3).Value = "14_ABASTO"
$ws.Cells.Item(82, 4).Value = 77
$ws.Cells.Item(83, 1).Value = "07:12:47"
$ws.Cells.Item(83, 2).Value = "08:33"
$ws.Cells.Item(83, 3).Value = "215C_EL PATO"
$ws.Cells.Item(83, 4).Value = 81
$ws.Cells.Item(84, 1).Value = "06:35:33"
$ws.Cells.Item(84, 2).Value = "08:34"
$ws.Cells.Item(84, 3).Value = "215C_EL PATO"
$ws.Cells.Item(84, 4).Value = 119
$ws.Cells.Item(85, 1).Value = "07:12:47"
$ws.Cells.Item(85, 2).Value = "08:47"
$ws.Cells.Item(85, 3).Value = "215A_EL PATO"
$ws.Cells.Item(85, 4).Value = 95
$ws.Cells.Item(85, 5).Value = "LP1912"
$ws.Cells.Item(86, 1).Value = "07:12:47"
$ws.Cells.Item(86, 2).Value = "08:51"
$ws.Cells.Item(86, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(86, 4).Value = 99
$ws.Cells.Item(86, 5).Value = "LP1912"
$ws.Cells.Item(87, 1).Value = "07:12:47"
$ws.Cells.Item(87, 2).Value = "08:59"
$ws.Cells.Item(87, 3).Value = "215B_EL PATO"
$ws.Cells.Item(87, 4).Value = 107
$ws.Cells.Item(87, 5).Value = "LP1912"

# ---- LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 07:12:47"
$ws.Cells.Item(3, 1).Value = "Total filas: 16"
$ws.Cells.Item(16, 1).Value = "07:12:47"
$ws.Cells.Item(16, 4).Value = 8
$ws.Cells.Item(18, 1).Value = "07:12:47"
$ws.Cells.Item(18, 4).Value = 81
$ws.Cells.Item(20, 1).Value = "07:12:47"
$ws.Cells.Item(20, 4).Value = 95
$ws.Cells.Item(21, 1).Value = "07:12:47"
$ws.Cells.Item(21, 2).Value = "08:59"
$ws.Cells.Item(21, 3).Value = "215B_EL PATO"
$ws.Cells.Item(21, 4).Value = 107
$ws.Cells.Item(21, 5).Value = "LP1912"

# ---- 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 07:12:47"
$ws.Cells.Item(7, 1).Value = "07:12:47"
$ws.Cells.Item(7, 4).Value = 15
$ws.Cells.Item(8, 1).Value = "07:12:47"
$ws.Cells.Item(8, 4).Value = 57
$ws.Cells.Item(10, 1).Value = "07:12:47"
$ws.Cells.Item(10, 4).Value = 70
